# feat: add 2022-Q4 data
#
# Plan (matches the target OOXML structure):
#  - The physical sheet that is currently "2022-Q3" (sheetId=2) is repurposed
#    in place to become "2022-Q4" and filled with the new quarter's fund
#    data (formatted like the "总计" sheet: same header/row styling and page
#    margins).
#  - A fresh copy of the original "2022-Q3" sheet (made before any edits) is
#    kept, unmodified, and renamed to "2022-Q3" - it ends up as a new
#    physical sheet (new sheetId) positioned after "2022-Q4".
#  - The "总计" summary sheet gets its existing "2022-Q3" label turned into
#    "2022-Q4", plus a new row re-adding a "2022-Q3" summary line.

$wb = $excel.ActiveWorkbook

$zj = $wb.Worksheets.Item("总计")
$origQ3 = $wb.Worksheets.Item("2022-Q3")

# --- 1. Preserve the original 2022-Q3 data by copying the sheet first.
#        The copy is placed right after the original and keeps its values
#        untouched - this becomes the new "2022-Q3" sheet. -------------------
$origQ3.Copy($null, $origQ3)
$q3Copy = $wb.ActiveSheet
$q3Copy.Name = "2022-Q3-new"

# --- 2. Turn the original sheet into "2022-Q4" and restyle it to match
#        "总计" (border/font style + page margins). ---------------------------
$origQ3.Name = "2022-Q4"

$zj.Range("B1").Copy()
$origQ3.Range("B1:H1").PasteSpecial(-4122)
$origQ3.Range("A2").PasteSpecial(-4122)

$zj.Range("B2").Copy()
$origQ3.Range("B2:H2").PasteSpecial(-4122)

$origQ3.PageSetup.LeftMargin = 54
$origQ3.PageSetup.RightMargin = 54
$origQ3.PageSetup.TopMargin = 72
$origQ3.PageSetup.BottomMargin = 72
$origQ3.PageSetup.HeaderMargin = 36
$origQ3.PageSetup.FooterMargin = 36

# --- helper: write a value as plain text (no leading '/NumberFormat residue
#             on the destination cell - paste-values-only keeps the target
#             cell's existing style untouched). -----------------------------
function Set-TextValue($ws, $addr, $val) {
    $helper = $ws.Range("Z100")
    $helper.NumberFormat = "@"
    $helper.Value = $val
    $helper.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $helper.Clear()
}

# --- 3. Populate "2022-Q4" with the new quarter's fund data. ----------------
$origQ3.Range("A2").Value = 0
Set-TextValue $origQ3 "B2" "090011"
$origQ3.Range("C2").Value = "大成核心双动力混合"
Set-TextValue $origQ3 "D2" "0.26"
Set-TextValue $origQ3 "E2" "93.43"
Set-TextValue $origQ3 "F2" "4.00"
Set-TextValue $origQ3 "G2" "0.0104"
$origQ3.Range("H2").Value = 4

# --- 4. Finish renaming the preserved copy back to "2022-Q3". ---------------
$q3Copy.Name = "2022-Q3"

# --- 5. Update the "总计" summary sheet: relabel the existing row as
#        2022-Q4 and add a fresh 2022-Q3 row underneath it. -----------------
$zj.Range("B2").Value = "2022-Q4"

$zj.Range("A2:D2").Copy()
$zj.Range("A3:D3").PasteSpecial(-4122)
$zj.Range("A3").Value = 1
$zj.Range("B3").Value = "2022-Q3"
$zj.Range("C3").Value = 1
$zj.Range("D3").Value = 0.01

# --- 6. Restore the original active tab (2022-Q3). ---------------------------
$wb.Worksheets.Item("2022-Q3").Activate()

Write-Output "done"
